$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(0.7652488218523669, 3.287053068188353, 1.287430089285444, 0.7652488218523665)
    3  = @(0.7521081881793897, 3.482316337766491, 1.367345623478399, 0.7352005490771096)
    4  = @(0.4688054747003463, 1.743119387866745, 0.2251811710404823, 0.4688054747003464)
    5  = @(0.4383441132901364, 1.841827598857237, 0.2412183778558042, 0.3904487462008647)
    6  = @(0.5529044951640609, 0.5015287438418512, 0.3922437331892626, 0.5529044951640609)
    7  = @(0.5293487356155337, 0.5249115170629832, 0.4138859651466079, 0.5022059794685393)
    8  = @(0.549793735567308,  230.9589342265355,  68.95711596322518,  0.5497937355673072)
    9  = @(0.5215083906738172, 244.9523546852538,  74.60437642749604,  0.4730347656325711)
    10 = @(0.6704659607956678, 5.057729332748186,  1.846388222750611,  0.6704659607956677)
    11 = @(0.6562443076201127, 5.286829133673952,  1.94719964562402,   0.6334989672791522)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}
